# #5: property boat&car done
# Add header row labels + legislator/meta columns (H:N) to the car ("汽車") sheet,
# matching the pattern already used on the land/building sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Row 1: header labels for columns B..N ----
# Give the new header cells (H1:N1) the same look (border/bold/etc.) as the
# existing header cells before filling them in.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial(-4122) | Out-Null

$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# ---- Row 2: fill in the new metadata columns (H2:N2) ----
# Give the new data cells (H2:N2) the same look as the existing data cells.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("H2:N2").PasteSpecial(-4122) | Out-Null

$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# Force the date to be stored as literal text (matches "date" column elsewhere
# in this workbook), not auto-converted to a date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2011-12-27"

$ws.Range("K2").Value = "楊瓊瓔"
$ws.Range("L2").Value = 854
$ws.Range("M2").Value = "tmpd1401"
$ws.Range("N2").Value = 44
